$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.443.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.78%  "
$ws.Range("D3").Value = "'1.770.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'306.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "'0.4297"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "'0.07199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'0.8481"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").Value = "'20.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "'1.781.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("D13").Value = "'6.432"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'5.234"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "'0.06907"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "'0.000008680"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "'26.442.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -12.09%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'11.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").Value = "'1.992.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").Value = "'152.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").Value = "'1.873"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.07%  "
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "'5.080"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "'114.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'1.744"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").Value = "'0.08962"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'0.7236"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("D33").Value = "'1.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'4.326"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").Value = "'2.750"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.87%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'1.079"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "'0.05156"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Value = "'0.01890"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'0.4924"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "'2.581"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.56%  "
$ws.Range("D43").Value = "'6.249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "'7.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("D45").Value = "'104.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'10.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "'0.4476"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'1.740"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.96%  "
